$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$values = @(
    'victoria secret workout',
    'victoria secret yoga capris',
    'victoria secret yoga pants',
    'victoria secret yoga pants for women',
    'victoria sport leggings',
    'victoria sport tight',
    'victorian fans for women',
    'victorian pants',
    'victorian style clothes for women',
    'victorias secret black leggings',
    'victorias secret capri',
    'victorias secret legging',
    'victorias secret leggings',
    'victorias secret leggings for women',
    'victorias secret leggins',
    'victorias secret sport leggings',
    'victorias secret tights',
    'violate leggings',
    'virus leggings women',
    'virus women compression pants',
    'virus womens compression',
    'visibility running waist reflective',
    'voluptuous women',
    'vs bed',
    'vs clothes for women',
    'vs leggings for women',
    'vs sport leggings',
    'w pants for women',
    'w womens clothing',
    'waist and hip compression',
    'waist compression',
    'waist compression for women',
    'waist compression women',
    'waist control tights',
    'waist high compression stockings',
    'waist high pantyhose',
    'waist hoop',
    'waist leggings',
    'waist length',
    'waist pack for running',
    'waist pack for women running',
    'waist pack reflective',
    'waist pack running 6',
    'waist running pack',
    'waist sinch',
    'waist support',
    'waist tight yoga leggings',
    'waist trainer capris',
    'waist training',
    'waist training for women',
    'waist training leggings',
    'waist warmer women',
    'waist workout',
    'waist yoga pants',
    'waistband fitness',
    'walking capris for women',
    'warm and form back brace',
    'warm compression pants boys',
    'warm compression pants mens',
    'warm compression pants women',
    'warm hose',
    'warm leggings for teens',
    'warm legings women',
    'warm long underwear women',
    'warm pants women',
    'warm running clothes for women',
    'warm running gear',
    'warm running gear women',
    'warm running pants women with pockets',
    'warm soccer gear',
    'warm trousers men',
    'warm tummy control leggings',
    'warm up pants for women',
    'warm up suits for women',
    'warm weather running gear women',
    'warm women leggings',
    'warm women tights',
    'warmup bands',
    'warmup pants for women',
    'warmup suits for women',
    'warmups women',
    'washing machine legs',
    'water capri',
    'water leggings women',
    'water pants women',
    'wear figs women',
    'wearable back support',
    'wearable flash drive',
    'wearable pain relief device',
    'wearable weight',
    'wears for women',
    'weather pants',
    'weather pocket chart',
    'weather sports',
    'weave knit wear',
    'web gear',
    'web knee brace',
    'web mail mil',
    'web med',
    'web thigh highs'
)
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
